# RedisCommands.xlsx - add GEOADD / GEOPOS / GEOHASH / GEODIST method mappings,
# mark a handful of "won't implement" rows as SKIP, flip a couple of finished
# flags, and refresh the hidden/filtered row state to match the autofilter
# (filters on column C == "FALSE").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows that are "not needed" / "not to be implemented" -> mark as SKIP ---
# (column C currently holds boolean FALSE; the new string value must be the
# very first new shared string written so it lands at the expected index)
$ws.Range("C19").Value = "SKIP"
foreach ($r in @(22, 23, 99, 100, 105, 108)) {
    $ws.Cells.Item($r, 3).Value = "SKIP"
}

# --- Flip a couple of "Finished" booleans from FALSE to TRUE ---
$ws.Range("C34").Value = $true

# --- New Geo commands: record the wrapper method name used to implement each ---
# (write in this exact order so new shared strings come out in the same order
# as the source commit)
$ws.Range("E139").Value = "GeoAddAsync"       # GEOADD
$ws.Range("E141").Value = "GeoPositionAsync"  # GEOPOS
$ws.Range("E140").Value = "GeoHashAsync"      # GEOHASH
$ws.Range("E142").Value = "GeoDistanceAsync"  # GEODIST

# GEOADD / GEOHASH / GEOPOS / GEODIST are now implemented
foreach ($r in @(139, 140, 141, 142)) {
    $ws.Cells.Item($r, 3).Value = $true
}

# --- Re-sync hidden rows with the autofilter (hides every row whose column C
# is no longer exactly FALSE, i.e. finished rows and the new SKIP rows) ---
$hiddenRows = @(19, 22, 23, 27, 42, 43, 82, 83, 85, 86, 89, 94, 95, 96, 97, 98, 99, 100, 101, 102, 103, 104, 105, 107, 108, 112)
foreach ($r in $hiddenRows) {
    $ws.Rows.Item($r).Hidden = $true
}

# --- Update the saved selection / scroll position ---
$ws.Activate() | Out-Null
$ws.Range("C142").Select() | Out-Null
